$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A22 must hold the literal text "09/23/2025" (matching the existing
# date cells in column A, which are stored as text, not real dates).
# Force text formatting before assigning the value so it isn't
# auto-converted to a date serial number, then clear the temporary
# formatting so the cell is left with the default style.
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "09/23/2025"
$ws.Range("A22").ClearFormats()

$ws.Range("B22").Value = 0.1311905007115779
$ws.Range("C22").Value = 0.8688094992884221
